$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    # Force text storage so numeric-looking strings are not converted to numbers
    $ws.Range($ref).NumberFormat = "@"
    $ws.Range($ref).Value = $val
    $ws.Range($ref).Style = "Normal"
}

Set-TextCell $ws "D2" "42.338.24"
$ws.Range("E2").Value = "  +1.34%  "
Set-TextCell $ws "D3" "2.272.47"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.03%  "
Set-TextCell $ws "D5" "307.01"
$ws.Range("E5").Value = "  +1.14%  "
Set-TextCell $ws "D6" "97.15"
$ws.Range("E6").Value = "  +5.70%  "
Set-TextCell $ws "D7" "0.529"
$ws.Range("E7").Value = "  -0.43%  "
Set-TextCell $ws "D9" "0.492"
$ws.Range("E9").Value = "  +1.93%  "
Set-TextCell $ws "D10" "35.73"
$ws.Range("E10").Value = "  +10.31%  "
Set-TextCell $ws "D11" "0.0795"
$ws.Range("E11").Value = "  -0.22%  "
$ws.Range("E12").Value = "  -2.20%  "
Set-TextCell $ws "D13" "6.66"
$ws.Range("E13").Value = "  -0.10%  "
Set-TextCell $ws "D14" "2.625.71"
$ws.Range("E14").Value = "  +0.32%  "
Set-TextCell $ws "D15" "14.41"
$ws.Range("E15").Value = "  +1.30%  "
Set-TextCell $ws "D16" "2.272.28"
$ws.Range("E16").Value = "  -0.53%  "
Set-TextCell $ws "D17" "0.791"
$ws.Range("E17").Value = "  +2.83%  "
Set-TextCell $ws "D18" "42.244.28"
$ws.Range("E18").Value = "  +1.34%  "
Set-TextCell $ws "D19" "12.50"
$ws.Range("E19").Value = "  +0.18%  "
Set-TextCell $ws "D20" "0.0₃0908"
$ws.Range("E20").Value = "  +0.34%  "
Set-TextCell $ws "D21" "5.96"
$ws.Range("E21").Value = "  +0.38%  "
Set-TextCell $ws "D22" "67.47"
$ws.Range("E22").Value = "  +0.52%  "
Set-TextCell $ws "D23" "240.53"
$ws.Range("E23").Value = "  +0.36%  "
$ws.Range("E24").Value = "  +0.47%  "
Set-TextCell $ws "D25" "1.93"
$ws.Range("E25").Value = "  +0.59%  "
Set-TextCell $ws "D26" "1.00"
$ws.Range("E26").Value = "  -0.10%  "
Set-TextCell $ws "D27" "23.74"
$ws.Range("E27").Value = "  -0.64%  "
Set-TextCell $ws "D28" "37.39"
$ws.Range("E28").Value = "  +6.87%  "
Set-TextCell $ws "D29" "9.47"
$ws.Range("E29").Value = "  -0.43%  "
Set-TextCell $ws "D30" "2.10"
$ws.Range("E30").Value = "  +1.19%  "
Set-TextCell $ws "D31" "158.89"
$ws.Range("E31").Value = "  -1.27%  "
Set-TextCell $ws "D32" "5.22"
$ws.Range("E32").Value = "  -0.70%  "
$ws.Range("E33").Value = "  +0.09%  "
Set-TextCell $ws "D34" "3.12"
Set-TextCell $ws "D35" "0.0741"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("B36").Value = "WEMIXToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextCell $ws "D36" "2.39"
$ws.Range("E36").Value = "  +1.13%  "
$ws.Range("B37").Value = "Celestia"
$ws.Range("C37").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
Set-TextCell $ws "D37" "16.86"
$ws.Range("E37").Value = "  -0.19%  "
Set-TextCell $ws "D38" "0.105"
$ws.Range("E38").Value = "  +1.15%  "
Set-TextCell $ws "D39" "1.83"
$ws.Range("E39").Value = "  +2.39%  "
$ws.Range("E40").Value = "  -1.45%  "
Set-TextCell $ws "D41" "4.09"
$ws.Range("E41").Value = "  +4.64%  "
Set-TextCell $ws "D42" "2.41"
$ws.Range("E42").Value = "  +13.70%  "
Set-TextCell $ws "D43" "1.996.78"
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell $ws "D44" "19.03"
$ws.Range("E44").Value = "  -0.95%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextCell $ws "D45" "0.0285"
$ws.Range("E45").Value = "  +2.11%  "
$ws.Range("E46").Value = "  +3.01%  "
Set-TextCell $ws "D47" "9.95"
$ws.Range("E47").Value = "  -3.40%  "
Set-TextCell $ws "D48" "52.84"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("E49").Value = "  +0.06%  "
Set-TextCell $ws "D50" "72.13"
Set-TextCell $ws "D51" "91.88"
$ws.Range("E51").Value = "  +1.11%  "
